# Adds a new game's stats for "LeBron James" as row 9 (pushing the
# previous "promedios" (averages) row down to row 10), and updates the
# "final" summary sheet so its LeBron James row pulls from the new
# averages row (row 10 instead of row 9).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("LeBron James")

# --- New game stats, currently sitting in the averages row (row 9) ---
# Capture them so nothing is lost, then overwrite row 9 with the new
# game's raw stats and push the averages down to row 10.

$ws.Range("A9:J9").ClearContents()

$ws.Range("A9").Value = 435
$ws.Range("B9").Value = 5
$ws.Range("C9").Value = 435
$ws.Range("D9").Value = 43
$ws.Range("E9").Value = 918
$ws.Range("F9").Value = 100
$ws.Range("G9").Value = 100.57
$ws.Range("H9").Value = 101.12
$ws.Range("I9").Value = 950

# --- Recreate the "promedios" (averages) row, now at row 10 ---
$ws.Range("A10").Formula = "=AVERAGE(A2:A9)"
$ws.Range("B10").Formula = "=AVERAGE(B2:B9)"
$ws.Range("C10").Formula = "=AVERAGE(C2:C9)"
$ws.Range("D10").Formula = "=AVERAGE(D2:D9)"
$ws.Range("E10").Formula = "=AVERAGE(E2:E9)"
$ws.Range("F10").Formula = "=AVERAGE(F2:F9)"
$ws.Range("G10").Formula = "=AVERAGE(G2:G9)"
$ws.Range("H10").Formula = "=AVERAGE(H2:H9)"
$ws.Range("I10").Formula = "=AVERAGE(I2:I9)"
$ws.Range("J10").Value = "promedios"

# --- Update the "final" summary sheet to pull from the new averages row ---
$final = $wb.Worksheets.Item("final")

$final.Range("B2").Formula = "='LeBron James'!A10"
$final.Range("C2").Formula = "='LeBron James'!B10"
$final.Range("D2").Formula = "='LeBron James'!C10"
$final.Range("E2").Formula = "='LeBron James'!D10"
$final.Range("F2").Formula = "='LeBron James'!E10"
$final.Range("G2").Formula = "='LeBron James'!F10"
$final.Range("H2").Formula = "='LeBron James'!G10"
$final.Range("I2").Formula = "='LeBron James'!H10"
$final.Range("J2").Formula = "='LeBron James'!I10"

$wb.Save()
